# Add a "cost-fix" (specific fixed cost) column to the Commodity and
# Process sheets of the rivus model workbook, matching the upstream
# rivus.py change that introduces a fixed (size-independent) maintenance
# cost parameter alongside the existing variable cost parameter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Commodity sheet: insert a new "cost-fix" column between
#    "cost-inv-var" (D) and "cost-var" (old E, now shifts to F).
# ---------------------------------------------------------------------
$wsCommodity = $wb.Worksheets.Item("Commodity")

$wsCommodity.Columns.Item(5).Insert() | Out-Null

$wsCommodity.Range("E1").Value = "cost-fix"
$wsCommodity.Range("E2").Value = 7
$wsCommodity.Range("E3").Value = 5
$wsCommodity.Range("E4").Value = 10
$wsCommodity.Range("E5").Value = 10
$wsCommodity.Range("E6").Formula = "=NA()"
$wsCommodity.Range("E7").Value = 0

# restore the selected cell on this sheet (it is not the active sheet)
$wsCommodity.Range("E8").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) Process sheet: insert a new "cost-fix" column between
#    "cost-inv-var" (C) and "cost-var" (old D, now shifts to E).
# ---------------------------------------------------------------------
$wsProcess = $wb.Worksheets.Item("Process")

$wsProcess.Columns.Item(4).Insert() | Out-Null

$wsProcess.Range("D1").Value = "cost-fix"
$wsProcess.Range("D2").Value = 0
$wsProcess.Range("D3").Value = 0
$wsProcess.Range("D4").Value = 0
$wsProcess.Range("D5").Value = 0
$wsProcess.Range("D6").Value = 0
$wsProcess.Range("D7").Value = 0
$wsProcess.Range("D8").Value = 0
$wsProcess.Range("D9").Value = 0

# matching input-message data validation, like the other header cells
$wsProcess.Range("D1").Validation.Add(0)
$wsProcess.Range("D1").Validation.InputTitle = "Specific fixed costs (€/kW)"
$wsProcess.Range("D1").Validation.InputMessage = "Size-dependent part for maintaining a plant."

# select + activate, making Process the active sheet/tab
$wsProcess.Range("B2").Select() | Out-Null
$wsProcess.Activate() | Out-Null
